$d = $word.ActiveDocument

$replacements = @(
    @("93-23=70", "94-90=4"),
    @("13+25=38", "23+75=98"),
    @("17+73=90", "99-37=62"),
    @("21+11=32", "47+37=84"),
    @("50-50=0", "33+46=79"),
    @("62-16=46", "6+6=12"),
    @("20+14=34", "81-63=18"),
    @("55+15=70", "1+93=94"),
    @("11+55=66", "54+6=60"),
    @("26+50=76", "90-16=74"),
    @("24-14=10", "42+55=97"),
    @("46-43=3", "73+3=76"),
    @("58-4=54", "44-29=15"),
    @("68+3=71", "77+13=90"),
    @("21-2=19", "50-30=20"),
    @("97-64=33", "86-20=66"),
    @("55+21=76", "30-18=12"),
    @("71+16=87", "30+40=70"),
    @("15+56=71", "24+68=92"),
    @("3+83=86", "2+6=8"),
    @("86-56=30", "36-20=16"),
    @("32+12=44", "46+22=68"),
    @("23+20=43", "4+73=77"),
    @("62-18=44", "90-15=75"),
    @("47+44=91", "75-42=33"),
    @("37+60=97", "30+52=82"),
    @("17+2=19", "39-39=0"),
    @("64-49=15", "95-11=84"),
    @("67-18=49", "65+1=66"),
    @("90-71=19", "46+46=92"),
    @("30+31=61", "21+48=69"),
    @("57+30=87", "87-73=14"),
    @("2+59=61", "46-25=21"),
    @("57-23=34", "58-29=29"),
    @("77-75=2", "98-64=34"),
    @("87-32=55", "45-16=29"),
    @("81-53=28", "63-31=32"),
    @("64+29=93", "69-17=52"),
    @("72+2=74", "63+7=70"),
    @("0+84=84", "86-66=20"),
    @("7+55=62", "2+80=82"),
    @("27+46=73", "92-14=78"),
    @("8+17=25", "98-49=49"),
    @("53-6=47", "86-74=12"),
    @("76-32=44", "44+1=45"),
    @("82+4=86", "11+56=67"),
    @("42-31=11", "37-7=30"),
    @("83-41=42", "65+26=91"),
    @("48-8=40", "25+40=65"),
    @("22+63=85", "31+1=32"),
    @("58+20=78", "66+27=93"),
    @("4+87=91", "24+7=31"),
    @("49-40=9", "24+66=90"),
    @("55+25=80", "49-4=45"),
    @("27-8=19", "77-76=1"),
    @("23+69=92", "24+2=26"),
    @("1+71=72", "62-19=43"),
    @("78-62=16", "10+6=16"),
    @("60-46=14", "19+79=98"),
    @("42-9=33", "35+31=66"),
    @("76+6=82", "13+47=60"),
    @("49+8=57", "92-32=60"),
    @("39+8=47", "15+44=59"),
    @("23-2=21", "49-32=17"),
    @("99-87=12", "71+20=91"),
    @("93-48=45", "58-12=46"),
    @("93-12=81", "38+49=87"),
    @("38+4=42", "24+10=34"),
    @("0+95=95", "86-44=42"),
    @("66-59=7", "71-61=10"),
    @("84-14=70", "79-70=9"),
    @("86-7=79", "21-1=20"),
    @("95-43=52", "70-12=58"),
    @("79-37=42", "17+43=60"),
    @("91-79=12", "16+80=96"),
    @("61-51=10", "66+8=74"),
    @("52+8=60", "53+22=75"),
    @("74-7=67", "99-5=94"),
    @("79-38=41", "3-0=3"),
    @("1+23=24", "5+33=38"),
    @("61-43=18", "91-47=44"),
    @("19+68=87", "95-23=72"),
    @("84-62=22", "50-18=32"),
    @("84+10=94", "85-20=65"),
    @("97-70=27", "49-36=13"),
    @("81-18=63", "60+33=93"),
    @("84-82=2", "4+80=84"),
    @("31-28=3", "21+8=29"),
    @("84-51=33", "15+76=91"),
    @("29-14=15", "45+23=68"),
    @("95-64=31", "25-6=19"),
    @("82-63=19", "7+80=87"),
    @("18-0=18", "43+40=83"),
    @("45+32=77", "84-73=11"),
    @("24+16=40", "73-58=15"),
    @("71+9=80", "55-28=27"),
    @("71-7=64", "97-2=95"),
    @("97-42=55", "73+26=99"),
    @("0+40=40", "67-53=14"),
    @("72-4=68", "57+24=81"),
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Write-Output "Done: applied $($replacements.Count) replacements"
